# Expand the single combined row of word-metrics into one row per word.
# Before: row 2 held "parsimony. colossal profit " as one entry.
# After:  each word gets its own row (2-7) with its own x/y/width/height values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2..7: text, x, y, width, height
$rows = @(
    @("parsimony. ", 654, 669,              98,  23),
    @("colossal ",   168, 719,              75,  23),
    @("profit ",     243, 719,              48,  23),
    @("masquerades ",775, 785.4,            122, 23),
    @("benevolent ", 587, 968.1999999999999,100, 23),
    @("bequeathed ", 292, 993.1999999999999,107, 23)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}
